$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'script.js'
$ws.Range('C2').Value = 'Une page d’accueil montrant (de manière dynamique) tous les articles disponibles à la vente.'
$ws.Range('D2').Value = 'Ouvrir sur la page d''accueil du site web dans un navigateur'
$ws.Range('E2').Value = 'Affichage de l''ensemble des produits'
$ws.Range('F2').Value = 'Affichage des différents canapé disponible / Pas de produit (erreur api)'

$ws.Range('B3').Value = 'product.js'
$ws.Range('C3').Value = 'Une page montrant un produit spécifique'
$ws.Range('D3').Value = 'Ouvrir une page produit d''un article'
$ws.Range('E3').Value = 'Affichage du canapé avec sa photo, son texte alternatif, son nom, son prix, sa description, les couleurs possibles'
$ws.Range('F3').Value = 'Affichage du canapé / Pas de produit, pas d''option, pas d''ajout possible (api)'

$ws.Range('B4').Value = 'product.js'
$ws.Range('C4').Value = 'Ajout du produit dans le panier'
$ws.Range('D4').Value = 'Au clic sur "ajout au panier", le produit est sauvegardé dans le localstorage du navigateur'
$ws.Range('E4').Value = 'Message confirmation d''ajout au panier et sauvegarde de l''item dans le localstorage'
$ws.Range('F4').Value = 'Message confirmation d''ajout au panier / Message d''erreur sur la quantité ou la couleur'

$ws.Range('B5').Value = 'product.js'
$ws.Range('C5').Value = 'Choix nombre de produit et de la couleur'
$ws.Range('D5').Value = 'Affichage de la couleur et du nombre de produit'
$ws.Range('E5').Value = 'Affichage du nombre d''article du produit et de la couleur selectionnée dans le menu déroulant'
$ws.Range('F5').Value = 'Affichage du nombre et de la couleur / Affichage en alerte si pas de couleur selectionnée ou pas le bon nombre de produit'

$ws.Range('B6').Value = 'cart.js'
$ws.Range('C6').Value = 'Une page panier dans laquelle s''affiche les produits choisis'
$ws.Range('D6').Value = 'Ouvrir la page panier du site web dans le navigateur'
$ws.Range('E6').Value = 'Affichage de l''ensemble des produits du localstorage avec le bon nombre et la bonne couleur'
$ws.Range('F6').Value = 'Si articles alors affichage de ceux ci / Affichage d''un message qui renvois à l''accueil '

$ws.Range('B7').Value = 'cart.js'
$ws.Range('C7').Value = 'Affichage du nombre et du prix total de la commande'
$ws.Range('D7').Value = 'Sur la page panier, le total affiche le nombre de produit et le prix total de la commande'
$ws.Range('E7').Value = 'Addition des différents aticles du panier et de leur prix'
$ws.Range('F7').Value = 'Prix total et nombre de l''ensemble des articles / Message qui renvois à l''accueil si le panier est vide'

$ws.Range('B8').Value = 'cart.js'
$ws.Range('C8').Value = 'Suppression d''un produit'
$ws.Range('D8').Value = 'Clique sur le bouton supprimer dans un produit du panier'
$ws.Range('E8').Value = 'Suppression de l''article complet, rechargement de la page et du prix et nombre d''article'
$ws.Range('F8').Value = 'Produit supprimer de la page panier / Message qui renvois à l''accueil si plus de produit dans le panier'

$ws.Range('B9').Value = 'cart.js'
$ws.Range('C9').Value = 'Modification du nombre de produit'
$ws.Range('D9').Value = 'Clique sur la quantité d''un produit du panier ( plus ou moins )'
$ws.Range('E9').Value = 'Modification de la quantité d''un article, rechargement de la page et du prix et du nombre d''article dans le total.'
$ws.Range('F9').Value = 'Modification du nombre sur l''article et donc du prix total et du nombre de produit dans le panier / Message d''erreur si pas le bon nombre'

$ws.Range('B10').Value = 'cart.js'
$ws.Range('C10').Value = 'Vérification formulaire '
$ws.Range('D10').Value = 'Entrer de valeurs dans les champs du formulaire par les utilisateurs'
$ws.Range('E10').Value = 'Formulaire completer avec les informations clients'
$ws.Range('F10').Value = 'Formulaire prêt à être envoyer / Erreur sur les champs non ou mal remplis'

$ws.Range('B11').Value = 'cart.js'
$ws.Range('C11').Value = 'Envoie de commande'
$ws.Range('D11').Value = 'Clique sur le bouton commander'
$ws.Range('E11').Value = 'Redirection vers la page confirmation de commande'
$ws.Range('F11').Value = 'Redirection vers la page confirmation de commande / Message d''erreur du formulaire et impossibilité d''envoyer la commande'

$ws.Range('B12').Value = 'confirmation.js'
$ws.Range('C12').Value = 'Confirmation de commande'
$ws.Range('D12').Value = 'Ouverture de la page confirmation dans le navigateur'
$ws.Range('E12').Value = 'Affichage du numéro de commande'
$ws.Range('F12').Value = 'Numéro de commande qui s''affiche / Pas de numéro de commande'

$ws.Range("G11").Select()
